$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" column (column R) is appended after the existing "2019"
# column (column Q). The header row (row 4) gets the year label, and
# each data row (5-14) gets its corresponding figure. Number/border
# formatting for the new column is copied from column Q so every new
# cell picks up the same style its Q-column neighbor already has.
$values = [ordered]@{
    4  = 2020
    5  = 2.1
    6  = 2.4
    7  = 1.4
    8  = 3.2
    9  = 2.4
    10 = 0.8
    11 = 2.2000000000000002
    12 = 4.5
    13 = 1.4
    14 = 3.2
}

foreach ($row in $values.Keys) {
    $qCell = $ws.Cells.Item($row, 17)   # column Q
    $rCell = $ws.Cells.Item($row, 18)   # column R

    # Clone Q's formatting onto R, then write the new figure.
    $qCell.Copy()
    $rCell.PasteSpecial(-4122)          # xlPasteFormats
    $rCell.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Move/extend the visible selection as recorded for this sheet.
$ws.Range("R16:R17").Select()
